$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 465-491 (dates as Excel serial numbers, B=nuovi pos., C=somma mobile 7gg., D=per 100mila abitanti)
$data = @(
    @(44539, 1, 2, 291.5451895043732),
    @(44540, 0, 2, 291.5451895043732),
    @(44541, 0, 2, 291.5451895043732),
    @(44542, 0, 1, 145.7725947521866),
    @(44543, 0, 1, 145.7725947521866),
    @(44544, 1, 2, 291.5451895043732),
    @(44545, 0, 2, 291.5451895043732),
    @(44546, 0, 1, 145.7725947521866),
    @(44547, 0, 1, 145.7725947521866),
    @(44548, 0, 1, 145.7725947521866),
    @(44550, 1, 2, 291.5451895043732),
    @(44551, 0, 2, 291.5451895043732),
    @(44552, 0, 1, 145.7725947521866),
    @(44553, 0, 1, 145.7725947521866),
    @(44554, 0, 1, 145.7725947521866),
    @(44555, 0, 1, 145.7725947521866),
    @(44556, 0, 1, 145.7725947521866),
    @(44557, 1, 1, 145.7725947521866),
    @(44558, 2, 3, 437.3177842565598),
    @(44559, 1, 4, 583.0903790087464),
    @(44560, 1, 5, 728.862973760933),
    @(44561, 0, 5, 728.862973760933),
    @(44562, 3, 8, 1166.180758017493),
    @(44563, 1, 9, 1311.953352769679),
    @(44564, 0, 8, 1166.180758017493),
    @(44565, 0, 6, 874.6355685131196),
    @(44566, 4, 9, 1311.953352769679),
)

$startRow = 465
$formatSourceRow = 464

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the date-column style (border/center/bold/date numfmt) from the last existing row
    $ws.Cells.Item($formatSourceRow, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "Added rows 465-491. New used range: $($ws.UsedRange.Address())"
